$d = $word.ActiveDocument

# Locate the original heading run "[NAME] Analysis:" without touching the
# paragraph mark (so the paragraph's own properties - spacing, rPr, rsids,
# paraId/textId, bookmark, etc. - are left completely untouched).
$rng = $d.Content
$found = $rng.Find.Execute("[NAME] Analysis:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Replace just the run text with two runs that share identical run
    # formatting (same as the original run) but split the wording into
    # "Baruch Consolidation " and "Analysis:" as two separate <w:r> elements.
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:noProof w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Baruch Consolidation </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:noProof w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>Analysis:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}
